$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old Address/Position columns entirely (cells + their formatting)
$ws.Range("B1:C3").Clear()

# Header cell (A1 already carries the bold/filled/bordered "header" look -
# overwriting the value alone does not disturb that formatting)
$ws.Range("A1").Value = "Serial"

# Data rows: WDPE059A50 .. WDPE059A80 (31 values) in rows 2..32
# (A2 already carries the bordered/centered "data" look from the old sheet)
for ($i = 50; $i -le 80; $i++) {
    $row = $i - 50 + 2
    $ws.Cells.Item($row, 1).Value = "WDPE059A$i"
}

# Stamp the remaining new rows (A3:A32) with the same look A2 already has
$ws.Range("A2").Copy()
$ws.Range("A3:A32").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Selection as saved in the workbook
$ws.Range("B21").Select()

Write-Host "done"
